# Applies updated crypto market data (price + 1h volume change) per the
# Sat Sep 14 11:38:55 UTC 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '59.769.04'
$ws.Range("E2").Value = '  +2.70%  '

# Row 3
$ws.Range("D3").Value = '2.416.38'
$ws.Range("E3").Value = '  +2.09%  '

# Row 4
$ws.Range("E4").Value = '  +0.06%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '550.64'
$ws.Range("E5").Value = '  +0.42%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.89'
$ws.Range("E6").Value = '  +2.57%  '

# Row 7
$ws.Range("E7").Value = '  -0.03%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.588'
$ws.Range("E8").Value = '  +3.94%  '

# Row 9
$ws.Range("E9").Value = '  -0.29%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.73'
$ws.Range("E10").Value = '  +2.92%  '

# Row 11
$ws.Range("E11").Value = '  -2.08%  '

# Row 12
$ws.Range("E12").Value = '  +0.02%  '

# Row 13
$ws.Range("E13").Value = '  +2.34%  '

# Row 14
$ws.Range("D14").Value = '2.847.08'
$ws.Range("E14").Value = '  +2.08%  '

# Row 15
$ws.Range("D15").Value = '59.755.81'
$ws.Range("E15").Value = '  +2.83%  '

# Row 16
$ws.Range("E16").Value = '  +0.39%  '

# Row 17
$ws.Range("D17").Value = '2.422.70'
$ws.Range("E17").Value = '  +2.37%  '

# Row 18
$ws.Range("E18").Value = '  +2.65%  '

# Row 19
$ws.Range("E19").Value = '  +0.74%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '329.54'
$ws.Range("E20").Value = '  -0.40%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.69'
$ws.Range("E21").Value = '  -2.98%  '

# Row 22
$ws.Range("E22").Value = '  -0.05%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.71'
$ws.Range("E23").Value = '  +3.49%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.172'
$ws.Range("E24").Value = '  +2.93%  '

# Row 25
$ws.Range("E25").Value = '  +3.92%  '

# Row 26
$ws.Range("E26").Value = '  -0.02%  '

# Row 27
$ws.Range("E27").Value = '  +1.11%  '

# Row 28
$ws.Range("D28").Value = '0.0₃0775'
$ws.Range("E28").Value = '  +4.31%  '

# Row 29
$ws.Range("E29").Value = '  +0.22%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '170.42'
$ws.Range("E30").Value = '  +0.04%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.15'
$ws.Range("E31").Value = '  -0.06%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.61'
$ws.Range("E32").Value = '  +0.95%  '

# Row 33
$ws.Range("E33").Value = '  +1.65%  '

# Row 35
$ws.Range("B35").Value = 'FirstDigitalUSD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.12%  '

# Row 36
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.28'
$ws.Range("E36").Value = '  +3.17%  '

# Row 37
$ws.Range("E37").Value = '  -0.29%  '

# Row 38
$ws.Range("E38").Value = '  +0.06%  '

# Row 39
$ws.Range("E39").Value = '  +0.41%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '314.52'
$ws.Range("E40").Value = '  +9.05%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.409'
$ws.Range("E41").Value = '  -0.65%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.65'
$ws.Range("E42").Value = '  -1.11%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '138.26'
$ws.Range("E43").Value = '  -3.40%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0967'
$ws.Range("E44").Value = '  +1.71%  '

# Row 45
$ws.Range("E45").Value = '  -0.10%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '19.38'
$ws.Range("E46").Value = '  +2.58%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.579'
$ws.Range("E47").Value = '  +2.25%  '

# Row 48
$ws.Range("E48").Value = '  +0.44%  '

# Row 49
$ws.Range("E49").Value = '  +0.25%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.53'
$ws.Range("E50").Value = '  -0.32%  '

# Row 51
$ws.Range("E51").Value = '  -0.33%  '
